# Regenerate merged AHB files
# - rename "_old" header suffixes to "_FV2410"
# - rename "_new" header suffixes to "_FV2504"
# - turn the A1:U64 range into a native Excel Table ("Table1")
# - freeze the header row (split below row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the "_old" -> "_FV2410" headers (columns A..J)
$ws.Range("A1").Value = "Segmentname_FV2410"
$ws.Range("B1").Value = "Segmentgruppe_FV2410"
$ws.Range("C1").Value = "Segment_FV2410"
$ws.Range("D1").Value = "Datenelement_FV2410"
$ws.Range("E1").Value = "Segment ID_FV2410"
$ws.Range("F1").Value = "Code_FV2410"
$ws.Range("G1").Value = "Qualifier_FV2410"
$ws.Range("H1").Value = "Beschreibung_FV2410"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2410"
$ws.Range("J1").Value = "Bedingung_FV2410"

# column K stays "diff"

# 2) Rename the "_new" -> "_FV2504" headers (columns L..U)
$ws.Range("L1").Value = "Segmentname_FV2504"
$ws.Range("M1").Value = "Segmentgruppe_FV2504"
$ws.Range("N1").Value = "Segment_FV2504"
$ws.Range("O1").Value = "Datenelement_FV2504"
$ws.Range("P1").Value = "Segment ID_FV2504"
$ws.Range("Q1").Value = "Code_FV2504"
$ws.Range("R1").Value = "Qualifier_FV2504"
$ws.Range("S1").Value = "Beschreibung_FV2504"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2504"
$ws.Range("U1").Value = "Bedingung_FV2504"

# 3) Convert the data range into a real Excel Table (adds xl/tables/table1.xml,
#    the worksheet <tableParts> element and the autoFilter on A1:U64).
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U64"), $null, 1)
$lo.Name = "Table1"

# 4) Freeze the header row (pane split under row 1, top-left of scrolling
#    area at A2) and keep the active selection on the lower-left pane.
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$null
